$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from E1 to F1 and set header text
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Cells.Item(1, 6).Value = "time_taken"

# Fill time_taken values for data rows 2-58
$ws.Cells.Item(2, 6).Value = "2021-10-05 10:52:27.889353"
$ws.Cells.Item(3, 6).Value = "2021-10-05 10:52:27.889367"
$ws.Cells.Item(4, 6).Value = "2021-10-05 10:52:27.889370"
$ws.Cells.Item(5, 6).Value = "2021-10-05 10:52:27.889374"
$ws.Cells.Item(6, 6).Value = "2021-10-05 10:52:27.889377"
$ws.Cells.Item(7, 6).Value = "2021-10-05 10:52:27.889380"
$ws.Cells.Item(8, 6).Value = "2021-10-05 10:52:27.889384"
$ws.Cells.Item(9, 6).Value = "2021-10-05 10:52:27.889387"
$ws.Cells.Item(10, 6).Value = "2021-10-05 10:52:27.889390"
$ws.Cells.Item(11, 6).Value = "2021-10-05 10:52:27.889393"
$ws.Cells.Item(12, 6).Value = "2021-10-05 10:52:27.889396"
$ws.Cells.Item(13, 6).Value = "2021-10-05 10:52:27.889399"
$ws.Cells.Item(14, 6).Value = "2021-10-05 10:52:27.889402"
$ws.Cells.Item(15, 6).Value = "2021-10-05 10:52:27.889405"
$ws.Cells.Item(16, 6).Value = "2021-10-05 10:52:27.889408"
$ws.Cells.Item(17, 6).Value = "2021-10-05 10:52:27.889411"
$ws.Cells.Item(18, 6).Value = "2021-10-05 10:52:27.889414"
$ws.Cells.Item(19, 6).Value = "2021-10-05 10:52:27.889417"
$ws.Cells.Item(20, 6).Value = "2021-10-05 10:52:27.889420"
$ws.Cells.Item(21, 6).Value = "2021-10-05 10:52:27.889423"
$ws.Cells.Item(22, 6).Value = "2021-10-05 10:52:27.889426"
$ws.Cells.Item(23, 6).Value = "2021-10-05 10:52:27.889429"
$ws.Cells.Item(24, 6).Value = "2021-10-05 10:52:27.889432"
$ws.Cells.Item(25, 6).Value = "2021-10-05 10:52:27.889435"
$ws.Cells.Item(26, 6).Value = "2021-10-05 10:52:27.889439"
$ws.Cells.Item(27, 6).Value = "2021-10-05 10:52:27.889442"
$ws.Cells.Item(28, 6).Value = "2021-10-05 10:52:27.889445"
$ws.Cells.Item(29, 6).Value = "2021-10-05 10:52:27.889448"
$ws.Cells.Item(30, 6).Value = "2021-10-05 10:52:27.889451"
$ws.Cells.Item(31, 6).Value = "2021-10-05 10:52:27.889454"
$ws.Cells.Item(32, 6).Value = "2021-10-05 10:52:27.889457"
$ws.Cells.Item(33, 6).Value = "2021-10-05 10:52:27.889460"
$ws.Cells.Item(34, 6).Value = "2021-10-05 10:52:27.889464"
$ws.Cells.Item(35, 6).Value = "2021-10-05 10:52:27.889467"
$ws.Cells.Item(36, 6).Value = "2021-10-05 10:52:27.889470"
$ws.Cells.Item(37, 6).Value = "2021-10-05 10:52:27.889473"
$ws.Cells.Item(38, 6).Value = "2021-10-05 10:52:27.889477"
$ws.Cells.Item(39, 6).Value = "2021-10-05 10:52:27.889480"
$ws.Cells.Item(40, 6).Value = "2021-10-05 10:52:27.889483"
$ws.Cells.Item(41, 6).Value = "2021-10-05 10:52:27.889486"
$ws.Cells.Item(42, 6).Value = "2021-10-05 10:52:27.889489"
$ws.Cells.Item(43, 6).Value = "2021-10-05 10:52:27.889493"
$ws.Cells.Item(44, 6).Value = "2021-10-05 10:52:27.889497"
$ws.Cells.Item(45, 6).Value = "2021-10-05 10:52:27.889500"
$ws.Cells.Item(46, 6).Value = "2021-10-05 10:52:27.889503"
$ws.Cells.Item(47, 6).Value = "2021-10-05 10:52:27.889506"
$ws.Cells.Item(48, 6).Value = "2021-10-05 10:52:27.889509"
$ws.Cells.Item(49, 6).Value = "2021-10-05 10:52:27.889512"
$ws.Cells.Item(50, 6).Value = "2021-10-05 10:52:27.889515"
$ws.Cells.Item(51, 6).Value = "2021-10-05 10:52:27.889518"
$ws.Cells.Item(52, 6).Value = "2021-10-05 10:52:27.889521"
$ws.Cells.Item(53, 6).Value = "2021-10-05 10:52:27.889525"
$ws.Cells.Item(54, 6).Value = "2021-10-05 10:52:27.889528"
$ws.Cells.Item(55, 6).Value = "2021-10-05 10:52:27.889531"
$ws.Cells.Item(56, 6).Value = "2021-10-05 10:52:27.889534"
$ws.Cells.Item(57, 6).Value = "2021-10-05 10:52:27.889537"
$ws.Cells.Item(58, 6).Value = "2021-10-05 10:52:27.889540"

$excel.CutCopyMode = 0

